$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '41.962.80'
$ws.Cells.Item(2, 5).Value = '  +5.99%  '
$ws.Cells.Item(3, 4).Value = '2.228.74'
$ws.Cells.Item(3, 5).Value = '  +2.82%  '
$ws.Cells.Item(4, 5).Value = '  +0.10%  '
$ws.Cells.Item(5, 5).Value = '  +1.78%  '
$ws.Cells.Item(6, 4).Value = '0.624'
$ws.Cells.Item(6, 5).Value = '  +0.47%  '
$ws.Cells.Item(7, 4).Value = '61.09'
$ws.Cells.Item(7, 5).Value = '  -2.73%  '
$ws.Cells.Item(8, 5).Value = '  +0.09%  '
$ws.Cells.Item(9, 4).Value = '0.401'
$ws.Cells.Item(9, 5).Value = '  +2.98%  '
$ws.Cells.Item(10, 4).Value = '58.97'
$ws.Cells.Item(10, 5).Value = '  +1.02%  '
$ws.Cells.Item(11, 5).Value = '  +5.22%  '
$ws.Cells.Item(12, 5).Value = '  -0.25%  '
$ws.Cells.Item(13, 4).Value = '2.559.59'
$ws.Cells.Item(13, 5).Value = '  +2.87%  '
$ws.Cells.Item(14, 4).Value = '15.64'
$ws.Cells.Item(14, 5).Value = '  -1.38%  '
$ws.Cells.Item(15, 4).Value = '21.74'
$ws.Cells.Item(15, 5).Value = '  +0.21%  '
$ws.Cells.Item(16, 4).Value = '0.798'
$ws.Cells.Item(16, 5).Value = '  -1.03%  '
$ws.Cells.Item(17, 5).Value = '  +2.01%  '
$ws.Cells.Item(18, 4).Value = '2.243.46'
$ws.Cells.Item(18, 5).Value = '  +3.66%  '
$ws.Cells.Item(19, 4).Value = '41.794.21'
$ws.Cells.Item(19, 5).Value = '  +5.53%  '
$ws.Cells.Item(20, 4).Value = '72.52'
$ws.Cells.Item(20, 5).Value = '  +1.23%  '
$ws.Cells.Item(21, 5).Value = '  +0.29%  '
$ws.Cells.Item(22, 4).Value = '6.02'
$ws.Cells.Item(22, 5).Value = '  +0.34%  '
$ws.Cells.Item(23, 4).Value = '250.08'
$ws.Cells.Item(23, 5).Value = '  +9.77%  '
$ws.Cells.Item(24, 4).Value = '0.999'
$ws.Cells.Item(24, 5).Value = '  -0.14%  '
$ws.Cells.Item(25, 5).Value = '  +1.45%  '
$ws.Cells.Item(26, 5).Value = '  -4.45%  '
$ws.Cells.Item(27, 4).Value = '9.69'
$ws.Cells.Item(27, 5).Value = '  +2.63%  '
$ws.Cells.Item(28, 5).Value = '  +3.49%  '
$ws.Cells.Item(29, 4).Value = '167.41'
$ws.Cells.Item(30, 5).Value = '  +1.34%  '
$ws.Cells.Item(31, 5).Value = '  -2.61%  '
$ws.Cells.Item(32, 4).Value = '2.64'
$ws.Cells.Item(32, 5).Value = '  -1.57%  '
$ws.Cells.Item(33, 5).Value = '  -0.33%  '
$ws.Cells.Item(34, 4).Value = '4.94'
$ws.Cells.Item(34, 5).Value = '  +5.06%  '
$ws.Cells.Item(35, 5).Value = '  +3.15%  '
$ws.Cells.Item(36, 5).Value = '  +1.99%  '
$ws.Cells.Item(37, 4).Value = '6.63'
$ws.Cells.Item(37, 5).Value = '  -4.80%  '
$ws.Cells.Item(38, 4).Value = '3.68'
$ws.Cells.Item(38, 5).Value = '  -2.77%  '
$ws.Cells.Item(39, 4).Value = '2.35'
$ws.Cells.Item(39, 5).Value = '  -1.52%  '
$ws.Cells.Item(40, 4).Value = '0.000257'
$ws.Cells.Item(41, 5).Value = '  +0.06%  '
$ws.Cells.Item(42, 5).Value = '  +5.57%  '
$ws.Cells.Item(43, 5).Value = '  -1.91%  '
$ws.Cells.Item(44, 5).Value = '  +8.64%  '
$ws.Cells.Item(45, 5).Value = '  +7.30%  '
$ws.Cells.Item(46, 4).Value = '98.78'
$ws.Cells.Item(46, 5).Value = '  -3.41%  '
$ws.Cells.Item(47, 5).Value = '  +0.54%  '
$ws.Cells.Item(48, 4).Value = '1.470.08'
$ws.Cells.Item(48, 5).Value = '  -2.83%  '
$ws.Cells.Item(49, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(49, 4).Value = '16.48'
$ws.Cells.Item(49, 5).Value = '  -7.17%  '
$ws.Cells.Item(50, 2).Value = 'HuobiToken'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(50, 4).Value = '2.81'
$ws.Cells.Item(50, 5).Value = '  +0.06%  '
$ws.Cells.Item(51, 5).Value = '  -1.31%  '
